$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 627, shifting the existing rows 627.. down by one
# (old row 627 becomes 628, ..., old row 668 becomes 669).
$ws.Rows(627).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row with the new data point (2026/01/11, 日, 13, 96).
# Column A holds a date-like string that must stay plain text (as every other
# row in this column does), so we briefly force a text number format while
# assigning the value and then restore the default "Normal" style so the new
# row doesn't pick up any extra cell formatting.
$ws.Cells.Item(627, 1).NumberFormat = "@"
$ws.Cells.Item(627, 1).Value = "2026/01/11"
$ws.Cells.Item(627, 1).Style = "Normal"

$ws.Cells.Item(627, 2).Value = "日"
$ws.Cells.Item(627, 3).Value = 13
$ws.Cells.Item(627, 4).Value = 96
